$d = $word.ActiveDocument

# Update the date line (first paragraph, outside the table)
$d.Content.Find.Execute("2025-02-28 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-01 Saturday", 2) | Out-Null

# Update each answer cell in the table by (row, col) position, since some
# old values repeat in the grid but map to different new values, a plain
# document-wide Find/Replace would not be safe here.
$t = $d.Tables(1)

$newValues = @(
    "63-7=56",
    "89-45=44",
    "74+6=80",
    "18+9=27",
    "77-30=47",
    "56+0=56",
    "36+56=92",
    "89-50=39",
    "29+39=68",
    "11+86=97",
    "73-40=33",
    "76+13=89",
    "79-66=13",
    "73-69=4",
    "48+3=51",
    "98-51=47",
    "56-0=56",
    "72+0=72",
    "20+4=24",
    "36+50=86",
    "91-83=8",
    "48-25=23",
    "46+13=59",
    "52-34=18",
    "85-16=69",
    "10+32=42",
    "36-19=17",
    "95-84=11",
    "56+33=89",
    "22+22=44",
    "7+4=11",
    "49+29=78",
    "49-35=14",
    "54-27=27",
    "1+39=40",
    "37+16=53",
    "70-49=21",
    "70-9=61",
    "68-55=13",
    "93-75=18",
    "1+23=24",
    "7+11=18",
    "12+14=26",
    "8+58=66",
    "51+15=66",
    "75-5=70",
    "1+1=2",
    "87-79=8",
    "35+45=80",
    "69-16=53",
    "18+57=75",
    "49+32=81",
    "23+43=66",
    "57-12=45",
    "27+25=52",
    "99-28=71",
    "93-4=89",
    "56+14=70",
    "15+27=42",
    "31-30=1",
    "66-35=31",
    "93+0=93",
    "97-48=49",
    "39-19=20",
    "20+50=70",
    "57+3=60",
    "88-54=34",
    "46+4=50",
    "2+77=79",
    "17+77=94",
    "90-28=62",
    "45-39=6",
    "83-13=70",
    "34-18=16",
    "24+58=82",
    "68-54=14",
    "99-91=8",
    "16+70=86",
    "15+78=93",
    "95-59=36",
    "20+50=70",
    "42-41=1",
    "72-20=52",
    "35-8=27",
    "32-21=11",
    "39+43=82",
    "36+54=90",
    "76-25=51",
    "91-44=47",
    "75-7=68",
    "8+51=59",
    "96-8=88",
    "25+31=56",
    "90+4=94",
    "18-10=8",
    "26+45=71",
    "95-8=87",
    "51-21=30",
    "64-2=62",
    "73-37=36"
)

$cols = 5
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $r = [int][math]::Floor($i / $cols) + 1
    $c = ($i % $cols) + 1
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $newValues[$i]
}

Write-Host "Done updating date and $($newValues.Length) table cells."
